$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "2024-09-02 18:07:55"
$ws.Range("K2").Value = "23h 59m 42s"

$ws.Range("C3").Value = "2024-09-02 18:07:55"
$ws.Range("K3").Value = "23h 59m 37s"

$ws.Range("C4").Value = "2023-09-10 18:08:25"

$ws.Range("C5").Value = "2023-09-04 18:08:25"

$ws.Range("C6").Value = "2024-10-07 18:08:25"

$ws.Range("C7").Value = "2023-09-10 18:08:25"
